$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("wenden","none","none"),
    @("geben","face/face030.jpg","face"),
    @("schicken","face/face018.jpg","face"),
    @("rufen","none","none"),
    @("leisten","face/face012.jpg","face"),
    @("bieten","face/face010.jpg","face"),
    @("enden","none","none"),
    @("quälen","face/face024.jpg","face"),
    @("heben","flower/flower028.jpg","flower"),
    @("opfern","none","none"),
    @("kommen","flower/flower003.jpg","flower"),
    @("fragen","flower/flower029.jpg","flower"),
    @("schalten","none","none"),
    @("wandern","face/face009.jpg","face"),
    @("helfen","face/face003.jpg","face"),
    @("klagen","none","none"),
    @("lehnen","flower/flower018.jpg","flower"),
    @("kehren","face/face016.jpg","face"),
    @("sparen","none","none"),
    @("schultern","flower/flower005.jpg","flower"),
    @("landen","face/face011.jpg","face"),
    @("hören","none","none"),
    @("rasen","flower/flower007.jpg","flower"),
    @("schweben","face/face026.jpg","face"),
    @("drohen","none","none"),
    @("tauschen","flower/flower020.jpg","flower"),
    @("fühlen","flower/flower010.jpg","flower"),
    @("orten","none","none"),
    @("herrschen","flower/flower031.jpg","flower"),
    @("drücken","flower/flower030.jpg","flower"),
    @("weigern","none","none"),
    @("biegen","face/face022.jpg","face"),
    @("zeugen","flower/flower025.jpg","flower"),
    @("ächzen","none","none"),
    @("dringen","flower/flower002.jpg","flower"),
    @("mühen","face/face023.jpg","face"),
    @("kosten","none","none"),
    @("heilen","face/face005.jpg","face"),
    @("spielen","face/face000.jpg","face"),
    @("dauern","none","none"),
    @("schreiben","face/face029.jpg","face"),
    @("tragen","flower/flower022.jpg","flower"),
    @("stören","none","none"),
    @("stopfen","flower/flower033.jpg","flower"),
    @("streichen","flower/flower015.jpg","flower"),
    @("bremsen","none","none"),
    @("achten","flower/flower012.jpg","flower"),
    @("tollen","face/face031.jpg","face")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
